$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITI")
$ws.Activate()

# Update column I (Status) for rows 2-33 to "Resolvido"
$range = $ws.Range("I2:I33")
$range.Value = "Resolvido"

# A couple of rows (I29:I30) previously carried a slightly different font;
# copy the formatting from a cell that already matches the rest of the
# column so the whole range ends up sharing one consistent style.
$ws.Range("I2").Copy()
$ws.Range("I29:I30").PasteSpecial(-4122)  # xlPasteFormats

# Highlight the whole updated range with a yellow fill
$range.Interior.Color = 65535

# Move the active selection to reflect where the user last worked
$ws.Range("I32:I33").Select()
